# Convert the M2Doc field-code paragraph (fldChar begin / instrText ... /
# fldChar end) into plain literal-text runs using "{" / "}" delimiters,
# as produced by the new TokenIteratorFieldRewriterSplit parser.

$d = $word.ActiveDocument

# Locate the paragraph that still contains the Word field (fldChar begin).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing the M2Doc field."
}

$r = $target.Range

$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"><w:r><w:t>{</w:t></w:r><w:r w:rsidR="00DE6D5A"><w:t>m</w:t></w:r><w:r w:rsidR="002033E1"><w:t>:</w:t></w:r><w:r w:rsidR="008B76C9"><w:t>'</w:t></w:r><w:r w:rsidR="00E806A4"><w:t>doc.html</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="008B76C9"><w:t>'.fromHTMLURI()</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p>
'@

$r.InsertXML($newParaXml)
